# Update the 13CThr_standards Concentrations sheet:
#  - remove the S2*/S3*/S4*/S6A rows (old rows 2-11)
#  - keep S6B, S6C, S7A, S7B, S7C (now rows 2-6)
#  - add new rows for S8A, S8B, S8C, S9A, S9B, S9C (rows 7-12)
#  - update selection / used range accordingly

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old data rows (below the header) entirely so the sheet
# only keeps the rows that remain in the final layout.
$ws.Rows("2:16").ClearContents()

# Data block: Sample, Threonine, Butyrate, Propionate, 2-hydroxybutyrate,
# 2-aminobutyrate, n-Propanol (columns A-G) -- header row (row 1) is
# unchanged.
$data = @(
    @("S6B", 12, 12, 12, 12, 12, 12),
    @("S6C", 12, 12, 12, 12, 12, 12),
    @("S7A", 15, 15, 15, 15, 15, 15),
    @("S7B", 15, 15, 15, 15, 15, 15),
    @("S7C", 15, 15, 15, 15, 15, 15),
    @("S8A",  9,  9,  9,  9,  9,  9),
    @("S8B",  9,  9,  9,  9,  9,  9),
    @("S8C",  9,  9,  9,  9,  9,  9),
    @("S9A", 15,  9,  9,  9,  9,  9),
    @("S9B", 15,  9,  9,  9,  9,  9),
    @("S9C", 15,  9,  9,  9,  9,  9)
)

$r = 2
foreach ($row in $data) {
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r++
}

# New used range is A1:G12; mirror the author's final UI selection (the
# whole of column F was selected, as if the column header was clicked).
$ws.Columns("F:F").Select()
